$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2025-08-30 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-31 Sunday", 2)

# Update the division problems in the single table, addressed by (row, col)
# so the chained old/new value overlaps (e.g. 39÷6=  -> 45÷9=, while a
# different cell's 45÷9= -> 50÷7=) do not clash with each other.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "45÷9="
$tbl.Cell(1, 2).Range.Text = "53÷4="
$tbl.Cell(1, 3).Range.Text = "63÷6="
$tbl.Cell(1, 4).Range.Text = "62÷9="
$tbl.Cell(1, 5).Range.Text = "46÷8="

$tbl.Cell(5, 1).Range.Text = "50÷7="
$tbl.Cell(5, 2).Range.Text = "89÷5="
$tbl.Cell(5, 3).Range.Text = "17÷4="
$tbl.Cell(5, 4).Range.Text = "74÷2="
$tbl.Cell(5, 5).Range.Text = "48÷5="

$tbl.Cell(9, 1).Range.Text = "19÷5="
$tbl.Cell(9, 2).Range.Text = "62÷2="
$tbl.Cell(9, 3).Range.Text = "65÷2="
$tbl.Cell(9, 4).Range.Text = "97÷9="
$tbl.Cell(9, 5).Range.Text = "80÷5="

$tbl.Cell(13, 1).Range.Text = "43÷4="
$tbl.Cell(13, 2).Range.Text = "56÷9="
$tbl.Cell(13, 3).Range.Text = "39÷9="
$tbl.Cell(13, 4).Range.Text = "25÷2="
$tbl.Cell(13, 5).Range.Text = "82÷3="

$tbl.Cell(17, 1).Range.Text = "94÷4="
$tbl.Cell(17, 2).Range.Text = "36÷6="
$tbl.Cell(17, 3).Range.Text = "10÷7="
$tbl.Cell(17, 4).Range.Text = "15÷4="
$tbl.Cell(17, 5).Range.Text = "44÷7="
